$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 227.53847
$ws.Range("I9").Value = 242.41667
$ws.Range("J9").Value = 49
$ws.Range("K9").Value = 242.41667
$ws.Range("L9").Value = 49
$ws.Range("M9").Value = -73.41667000000001
$ws.Range("N9").Value = -387
$ws.Range("H40").Value = 8486.117
$ws.Range("I40").Value = 7895.3335
$ws.Range("J40").Value = 8808.362999999999
$ws.Range("K40").Value = 7895.3335
$ws.Range("L40").Value = 8808.362999999999
$ws.Range("M40").Value = -7720.3335
$ws.Range("N40").Value = -9158.362999999999
$ws.Range("H112").Value = 3065.9473
$ws.Range("J112").Value = 3170.4243
$ws.Range("L112").Value = 9511.2729
$ws.Range("N112").Value = -11727.2729
$ws.Range("H131").Value = 4610.4707
$ws.Range("I131").Value = 1188.909
$ws.Range("K131").Value = 3566.727
$ws.Range("M131").Value = 1473.273
$ws.Range("H135").Value = 2826.6316
$ws.Range("I135").Value = 2087
$ws.Range("K135").Value = 18783
$ws.Range("M135").Value = -16248
$ws.Range("H138").Value = 2667.59
$ws.Range("J138").Value = 2816.2598
$ws.Range("L138").Value = 8448.779399999999
$ws.Range("N138").Value = -18728.7794
$ws.Range("H141").Value = 9258.5625
$ws.Range("I141").Value = 7971.3447
$ws.Range("J141").Value = 21701.666
$ws.Range("K141").Value = 23914.0341
$ws.Range("L141").Value = 65104.99800000001
$ws.Range("M141").Value = -18734.0341
$ws.Range("N141").Value = -75464.99800000001

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 22983344
$ws.Range("I2").Value = 40660990
$ws.Range("K2").Value = 40660990
$ws.Range("M2").Value = -40660877
$ws.Range("H55").Value = 34166.332
$ws.Range("J55").Value = 39799.6
$ws.Range("L55").Value = 39799.6
$ws.Range("N55").Value = -40429.6
$ws.Range("H63").Value = 3324.875
$ws.Range("I63").Value = 2282
$ws.Range("J63").Value = 3950.6
$ws.Range("K63").Value = 2282
$ws.Range("L63").Value = 3950.6
$ws.Range("M63").Value = -1596
$ws.Range("N63").Value = -5322.6
$ws.Range("H66").Value = 3324.875
$ws.Range("I66").Value = 2282
$ws.Range("J66").Value = 3950.6
$ws.Range("K66").Value = 11410
$ws.Range("L66").Value = 19753
$ws.Range("M66").Value = -7978
$ws.Range("N66").Value = -26617
$ws.Range("H88").Value = 1733.5
$ws.Range("J88").Value = 1103.2858
$ws.Range("L88").Value = 1103.2858
$ws.Range("N88").Value = -1915.2858
$ws.Range("H91").Value = 1733.5
$ws.Range("J91").Value = 1103.2858
$ws.Range("L91").Value = 1103.2858
$ws.Range("N91").Value = -3911.2858
$ws.Range("H116").Value = 22983344
$ws.Range("I116").Value = 40660990
$ws.Range("K116").Value = 40660990
$ws.Range("M116").Value = -40658696

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 22983344
$ws.Range("I3").Value = 40660990
$ws.Range("K3").Value = 40660990
$ws.Range("M3").Value = -40660876
$ws.Range("H20").Value = 2245.7673
$ws.Range("I20").Value = 2054.5806
$ws.Range("K20").Value = 2054.5806
$ws.Range("M20").Value = -1807.5806
$ws.Range("H105").Value = 2132.4736
$ws.Range("I105").Value = 1907.4375
$ws.Range("J105").Value = 3332.6667
$ws.Range("K105").Value = 1907.4375
$ws.Range("L105").Value = 3332.6667
$ws.Range("M105").Value = -160.4375
$ws.Range("N105").Value = -6826.6667

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2882.2124
$ws.Range("I31").Value = 2125.1833
$ws.Range("K31").Value = 2125.1833
$ws.Range("M31").Value = -1830.1833
$ws.Range("H34").Value = 2882.2124
$ws.Range("I34").Value = 2125.1833
$ws.Range("K34").Value = 2125.1833
$ws.Range("M34").Value = -1923.1833
$ws.Range("H58").Value = 1743.069
$ws.Range("I58").Value = 1502.25
$ws.Range("J58").Value = 2039.4615
$ws.Range("K58").Value = 1502.25
$ws.Range("L58").Value = 2039.4615
$ws.Range("M58").Value = -1299.25
$ws.Range("N58").Value = -2445.4615
$ws.Range("H62").Value = 333344260
$ws.Range("I62").Value = 500008900
$ws.Range("J62").Value = 15000
$ws.Range("K62").Value = 500008900
$ws.Range("L62").Value = 15000
$ws.Range("M62").Value = -500008276
$ws.Range("N62").Value = -16248
$ws.Range("H65").Value = 333344260
$ws.Range("I65").Value = 500008900
$ws.Range("J65").Value = 15000
$ws.Range("K65").Value = 2500044500
$ws.Range("L65").Value = 75000
$ws.Range("M65").Value = -2500041380
$ws.Range("N65").Value = -81240
$ws.Range("H68").Value = 51710.4
$ws.Range("J68").Value = 49749.25
$ws.Range("L68").Value = 49749.25
$ws.Range("N68").Value = -51247.25
$ws.Range("H71").Value = 51710.4
$ws.Range("J71").Value = 49749.25
$ws.Range("L71").Value = 149247.75
$ws.Range("N71").Value = -156735.75
$ws.Range("H99").Value = 12888.228
$ws.Range("J99").Value = 17172.908
$ws.Range("L99").Value = 17172.908
$ws.Range("N99").Value = -20168.908
$ws.Range("H126").Value = 12888.228
$ws.Range("J126").Value = 17172.908
$ws.Range("L126").Value = 51518.724
$ws.Range("N126").Value = -56458.724
$ws.Range("H132").Value = 6417.971
$ws.Range("I132").Value = 6220.6772
$ws.Range("J132").Value = 7947
$ws.Range("K132").Value = 18662.0316
$ws.Range("L132").Value = 23841
$ws.Range("M132").Value = -16132.0316
$ws.Range("N132").Value = -28901
$ws.Range("H134").Value = 3650.4167
$ws.Range("I134").Value = 3635.2173
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 10905.6519
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -8370.651899999999
$ws.Range("N134").Value = -17070
$ws.Range("H136").Value = 1743.069
$ws.Range("I136").Value = 1502.25
$ws.Range("J136").Value = 2039.4615
$ws.Range("K136").Value = 4506.75
$ws.Range("L136").Value = 6118.3845
$ws.Range("M136").Value = -1956.75
$ws.Range("N136").Value = -11218.3845

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 791
$ws.Range("I34").Value = 402.75
$ws.Range("K34").Value = 1208.25
$ws.Range("M34").Value = -1124.25
$ws.Range("H104").Value = 7637.1665
$ws.Range("I104").Value = 3331
$ws.Range("K104").Value = 9993
$ws.Range("M104").Value = -7372
$ws.Range("H113").Value = 2269.8235
$ws.Range("J113").Value = 1702.7142
$ws.Range("L113").Value = 5108.142599999999
$ws.Range("N113").Value = -9448.142599999999
$ws.Range("H122").Value = 925.2143
$ws.Range("I122").Value = 807.625
$ws.Range("J122").Value = 1082
$ws.Range("K122").Value = 7268.625
$ws.Range("L122").Value = 9738
$ws.Range("M122").Value = -4818.625
$ws.Range("N122").Value = -14638
$ws.Range("H129").Value = 959.5
$ws.Range("J129").Value = 1000
$ws.Range("L129").Value = 3000
$ws.Range("N129").Value = -13000
$ws.Range("H131").Value = 755350.25
$ws.Range("J131").Value = 1967.3
$ws.Range("L131").Value = 5901.9
$ws.Range("N131").Value = -15981.9

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7500
$ws.Range("J70").Value = 7500
$ws.Range("L70").Value = 7500
$ws.Range("N70").Value = -8040
$ws.Range("H73").Value = 7500
$ws.Range("J73").Value = 7500
$ws.Range("L73").Value = 7500
$ws.Range("N73").Value = -9372
$ws.Range("H80").Value = 10873570
$ws.Range("I80").Value = 20836268
$ws.Range("K80").Value = 20836268
$ws.Range("M80").Value = -20835270
$ws.Range("H83").Value = 10873570
$ws.Range("I83").Value = 20836268
$ws.Range("K83").Value = 104181340
$ws.Range("M83").Value = -104176348
$ws.Range("H113").Value = 4037.6924
$ws.Range("I113").Value = 3849.6667
$ws.Range("J113").Value = 4198.857
$ws.Range("K113").Value = 3849.6667
$ws.Range("L113").Value = 4198.857
$ws.Range("M113").Value = -1679.6667
$ws.Range("N113").Value = -8538.857
$ws.Range("H126").Value = 7206
$ws.Range("I126").Value = 4988
$ws.Range("K126").Value = 14964
$ws.Range("M126").Value = -12494
$ws.Range("H132").Value = 4778.9316
$ws.Range("I132").Value = 4704.606
$ws.Range("K132").Value = 14113.818
$ws.Range("M132").Value = -11583.818

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 894.125
$ws.Range("I55").Value = 834.5
$ws.Range("K55").Value = 834.5
$ws.Range("M55").Value = -661.5
$ws.Range("H68").Value = 23811616
$ws.Range("I68").Value = 27779636
$ws.Range("K68").Value = 27779636
$ws.Range("M68").Value = -27778887
$ws.Range("H71").Value = 23811616
$ws.Range("I71").Value = 27779636
$ws.Range("K71").Value = 138898180
$ws.Range("M71").Value = -138894436
$ws.Range("H132").Value = 17210.559
$ws.Range("I132").Value = 18951.03
$ws.Range("J132").Value = 11073.105
$ws.Range("K132").Value = 56853.09
$ws.Range("L132").Value = 33219.315
$ws.Range("M132").Value = -54323.09
$ws.Range("N132").Value = -38279.315

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 166665
$ws.Range("J141").Value = 166665
$ws.Range("L141").Value = 166665
$ws.Range("N141").Value = -177025
